$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: literal header-like values mirroring row 3 (R3:U3 = 2,3,4,5)
$ws.Range("R33").Value = 2
$ws.Range("S33").Value = 3
$ws.Range("T33").Value = 4
$ws.Range("U33").Value = 5

# Row 35: summary totals for columns R:U (row 34 left blank)
$ws.Range("R35").Formula = "=SUM(R4:R31)"
$ws.Range("S35:U35").Formula = "=SUM(S4:S31)"

# Move the active selection to reflect where the user ended up editing
$ws.Range("X30").Select()
